$d = $word.ActiveDocument

# --- Part 2: title "Cool Box" -> "Image Gallery" ---
$d.Content.Find.Execute("Cool Box", $false, $false, $false, $false, $false, $true, 1, $false, "Image Gallery", 2)

# --- Part 2: hyperlink URL (avoid exact full-run match so the Hyperlink rStyle survives) ---
$d.Content.Find.Execute("saadat176.github.io/baigs7/Assignment3/Pt2/Pt2-CoolBox.html", $false, $false, $false, $false, $false, $true, 1, $false, "saadat176.github.io/baigs7/Assignment4/Pt2/Pt2-ImageGallery.html", 2)

# --- Part 3: title "College Website" -> "Bouncing Balls" ---
$d.Content.Find.Execute("College Website", $false, $false, $false, $false, $false, $true, 1, $false, "Bouncing Balls", 2)

# --- Part 3: hyperlink URL ---
$d.Content.Find.Execute("saadat176.github.io/baigs7/Assignment3/Pt3/Pt3-CollegeWebsite.html", $false, $false, $false, $false, $false, $true, 1, $false, "saadat176.github.io/baigs7/Assignment4/Pt3/Pt3-BouncingBalls.html", 2)

# --- Part 4: title "My Cool Website" -> "Evil Circle vs Bouncing Balls" ---
$d.Content.Find.Execute("My Cool Website", $false, $false, $false, $false, $false, $true, 1, $false, "Evil Circle vs Bouncing Balls", 2)

# --- Part 4: hyperlink URL ---
$d.Content.Find.Execute("saadat176.github.io/baigs7/Assignment3/Pt4/Pt4-MyCoolWebsite.html", $false, $false, $false, $false, $false, $true, 1, $false, "saadat176.github.io/baigs7/Assignment4/Pt4/Pt4-BBvEvilCircle.html", 2)

# --- Remove one of the two trailing line breaks after the last (Part 4) hyperlink ---
$hlinkCount = $d.Hyperlinks.Count
$lastLink = $d.Hyperlinks.Item($hlinkCount)
$afterLink = $lastLink.Range.End
$brRange = $d.Range($afterLink, $afterLink + 1)
$brRange.Delete()

Write-Output "done"
